$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text changes ---
$ws.Range("A8").Value = "Volume 32   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/2/2025  Through  6/8/2025"

# --- Numeric cell changes ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 9
$ws.Range("I14").Value = 40
$ws.Range("J14").Value = 48
$ws.Range("K14").Value = -16.666666666666
$ws.Range("L14").Value = -24.528301886792
$ws.Range("M14").Value = -14.893617021276
$ws.Range("N14").Value = -80.487804878048

# Row 15
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = -14.285714285714
$ws.Range("F15").Value = 40
$ws.Range("G15").Value = 36
$ws.Range("H15").Value = 11.111111111111
$ws.Range("I15").Value = 230
$ws.Range("J15").Value = 187
$ws.Range("K15").Value = 22.994652406417
$ws.Range("L15").Value = 29.213483146067
$ws.Range("M15").Value = 88.524590163934
$ws.Range("N15").Value = -22.818791946308

# Row 16
$ws.Range("C16").Value = 104
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 398
$ws.Range("G16").Value = 378
$ws.Range("H16").Value = 5.291005291005
$ws.Range("I16").Value = 1969
$ws.Range("J16").Value = 2083
$ws.Range("K16").Value = -5.472875660105
$ws.Range("L16").Value = 2.392095683827
$ws.Range("M16").Value = 7.654455986878
$ws.Range("N16").Value = -71.438932404989

# Row 17
$ws.Range("C17").Value = 206
$ws.Range("D17").Value = 164
$ws.Range("E17").Value = 25.609756097561
$ws.Range("F17").Value = 717
$ws.Range("G17").Value = 713
$ws.Range("H17").Value = 0.561009817671
$ws.Range("I17").Value = 3604
$ws.Range("J17").Value = 3419
$ws.Range("K17").Value = 5.410938871014
$ws.Range("L17").Value = 9.544072948328
$ws.Range("M17").Value = 91.702127659574
$ws.Range("N17").Value = -2.065217391304

# Row 18
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 45
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 195
$ws.Range("G18").Value = 246
$ws.Range("H18").Value = -20.731707317073
$ws.Range("I18").Value = 1223
$ws.Range("J18").Value = 1261
$ws.Range("K18").Value = -3.013481363996
$ws.Range("L18").Value = -7.767722473604
$ws.Range("M18").Value = -9.407407407407
$ws.Range("N18").Value = -84.875092752906

# Row 19
$ws.Range("C19").Value = 182
$ws.Range("D19").Value = 183
$ws.Range("E19").Value = -0.546448087431
$ws.Range("F19").Value = 710
$ws.Range("G19").Value = 722
$ws.Range("H19").Value = -1.662049861495
$ws.Range("I19").Value = 3849
$ws.Range("J19").Value = 3913
$ws.Range("K19").Value = -1.635573728596
$ws.Range("L19").Value = 14.451382694023
$ws.Range("M19").Value = 107.829373650108
$ws.Range("N19").Value = 24.081237911025

# Row 20
$ws.Range("C20").Value = 79
$ws.Range("D20").Value = 82
$ws.Range("E20").Value = -3.658536585365
$ws.Range("F20").Value = 397
$ws.Range("G20").Value = 316
$ws.Range("H20").Value = 25.632911392405
$ws.Range("I20").Value = 1931
$ws.Range("J20").Value = 1732
$ws.Range("K20").Value = 11.489607390300
$ws.Range("L20").Value = -16.695427092321
$ws.Range("M20").Value = 125.320886814469
$ws.Range("N20").Value = -71.101466626758

# Row 21
$ws.Range("C21").Value = 620
$ws.Range("D21").Value = 589
$ws.Range("E21").Value = 5.263157894736
$ws.Range("F21").Value = 2463
$ws.Range("G21").Value = 2420
$ws.Range("H21").Value = 1.776859504132
$ws.Range("I21").Value = 12846
$ws.Range("J21").Value = 12643
$ws.Range("K21").Value = 1.605631574784
$ws.Range("L21").Value = 3.172435948919
$ws.Range("M21").Value = 61.849565326949
$ws.Range("N21").Value = -55.622344284381

# Row 22
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = 5
$ws.Range("I22").Value = 131
$ws.Range("J22").Value = 152
$ws.Range("K22").Value = -13.815789473684
$ws.Range("L22").Value = -1.503759398496
$ws.Range("M22").Value = -12.666666666666

# Row 23
$ws.Range("C23").Value = 22
$ws.Range("E23").Value = -24.137931034482
$ws.Range("F23").Value = 113
$ws.Range("H23").Value = -4.237288135593
$ws.Range("I23").Value = 660
$ws.Range("J23").Value = 726
$ws.Range("K23").Value = -9.090909090909
$ws.Range("L23").Value = -13.272010512483
$ws.Range("M23").Value = 54.205607476635

# Row 24
$ws.Range("C24").Value = 347
$ws.Range("D24").Value = 267
$ws.Range("E24").Value = 29.962546816479
$ws.Range("F24").Value = 1502
$ws.Range("G24").Value = 1107
$ws.Range("H24").Value = 35.682023486901
$ws.Range("I24").Value = 7689
$ws.Range("J24").Value = 7035
$ws.Range("K24").Value = 9.296375266524
$ws.Range("L24").Value = 0.654535934022
$ws.Range("M24").Value = 47.298850574712

# Row 25
$ws.Range("C25").Value = 113
$ws.Range("D25").Value = 99
$ws.Range("E25").Value = 14.141414141414
$ws.Range("F25").Value = 484
$ws.Range("G25").Value = 437
$ws.Range("H25").Value = 10.755148741418
$ws.Range("I25").Value = 2523
$ws.Range("J25").Value = 2834
$ws.Range("K25").Value = -10.973888496824
$ws.Range("L25").Value = -23.545454545454

# Row 26
$ws.Range("C26").Value = 248
$ws.Range("D26").Value = 236
$ws.Range("E26").Value = 5.084745762711
$ws.Range("F26").Value = 983
$ws.Range("G26").Value = 970
$ws.Range("H26").Value = 1.340206185567
$ws.Range("I26").Value = 4712
$ws.Range("J26").Value = 4691
$ws.Range("K26").Value = 0.447665742911
$ws.Range("L26").Value = 4.293935369632
$ws.Range("M26").Value = 0.170068027210

# Row 27
$ws.Range("C27").Value = 16
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = 6.666666666666
$ws.Range("F27").Value = 46
$ws.Range("G27").Value = 51
$ws.Range("H27").Value = -9.803921568627
$ws.Range("I27").Value = 289
$ws.Range("J27").Value = 291
$ws.Range("K27").Value = -0.687285223367
$ws.Range("L27").Value = -2.693602693602

# Row 28
$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 26
$ws.Range("E28").Value = -7.692307692307
$ws.Range("F28").Value = 100
$ws.Range("G28").Value = 105
$ws.Range("H28").Value = -4.761904761904
$ws.Range("I28").Value = 475
$ws.Range("J28").Value = 531
$ws.Range("K28").Value = -10.546139359698
$ws.Range("L28").Value = 3.036876355748

# Row 29
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 125
$ws.Range("F29").Value = 24
$ws.Range("G29").Value = 30
$ws.Range("H29").Value = -20
$ws.Range("I29").Value = 122
$ws.Range("J29").Value = 152
$ws.Range("K29").Value = -19.736842105263
$ws.Range("L29").Value = -11.594202898550
$ws.Range("M29").Value = -30.681818181818
$ws.Range("N29").Value = -76.806083650190

# Row 30
$ws.Range("C30").Value = 8
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 21
$ws.Range("G30").Value = 26
$ws.Range("H30").Value = -19.230769230769
$ws.Range("I30").Value = 106
$ws.Range("J30").Value = 124
$ws.Range("K30").Value = -14.516129032258
$ws.Range("L30").Value = -7.826086956521
$ws.Range("M30").Value = -28.378378378378
$ws.Range("N30").Value = -77.637130801687

# Row 31
$ws.Range("D31").Value = 2
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 3
$ws.Range("I31").Value = 9
$ws.Range("J31").Value = 14
$ws.Range("K31").Value = -35.714285714285
$ws.Range("L31").Value = -10

# Row 33
$ws.Range("D33").Value = 4
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = -83.333333333333
$ws.Range("J33").Value = 21
$ws.Range("K33").Value = -47.619047619047
